$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.689.96'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '2.100.17'
$ws.Range("E3").Value = '  +1.95%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'234.11"
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("D6").Value = "'0.622"
$ws.Range("E6").Value = '  +0.80%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = "'57.78"
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("D9").Value = "'0.390"
$ws.Range("E9").Value = '  +1.76%  '
$ws.Range("D10").Value = "'0.0776"
$ws.Range("E10").Value = '  +2.28%  '
$ws.Range("D12").Value = '2.412.60'
$ws.Range("E12").Value = '  +2.06%  '
$ws.Range("D13").Value = "'14.48"
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("D14").Value = "'21.41"
$ws.Range("E14").Value = '  +2.99%  '
$ws.Range("D15").Value = "'0.780"
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").Value = "'5.20"
$ws.Range("E16").Value = '  +1.03%  '
$ws.Range("D17").Value = '2.106.57'
$ws.Range("E17").Value = '  +2.22%  '
$ws.Range("D18").Value = '37.679.83'
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("D19").Value = "'6.18"
$ws.Range("E19").Value = '  -1.61%  '
$ws.Range("D20").Value = "'70.11"
$ws.Range("E20").Value = '  +1.25%  '
$ws.Range("D21").Value = '0.0₃0820'
$ws.Range("E21").Value = '  +1.05%  '
$ws.Range("D22").Value = "'226.93"
$ws.Range("E22").Value = '  +0.74%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("E24").Value = '  -1.76%  '
$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = '  +0.47%  '
$ws.Range("D26").Value = "'168.36"
$ws.Range("E26").Value = '  +1.50%  '
$ws.Range("D27").Value = "'8.91"
$ws.Range("E27").Value = '  +1.16%  '
$ws.Range("E28").Value = '  +1.50%  '
$ws.Range("E29").Value = '  -3.05%  '
$ws.Range("D30").Value = "'19.41"
$ws.Range("E30").Value = '  +1.85%  '
$ws.Range("E31").Value = '  +0.52%  '
$ws.Range("D32").Value = "'4.61"
$ws.Range("E32").Value = '  +2.66%  '
$ws.Range("E33").Value = '  +0.70%  '
$ws.Range("D34").Value = "'0.0620"
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("D35").Value = "'4.55"
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("D36").Value = "'3.44"
$ws.Range("E36").Value = '  +5.37%  '
$ws.Range("D37").Value = "'1.81"
$ws.Range("E37").Value = '  +4.20%  '
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("D39").Value = "'5.37"
$ws.Range("E39").Value = '  -8.62%  '
$ws.Range("D40").Value = "'0.0988"
$ws.Range("E40").Value = '  +6.17%  '
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("D42").Value = '1.472.41'
$ws.Range("E42").Value = '  -0.65%  '
$ws.Range("D43").Value = "'96.18"
$ws.Range("E43").Value = '  -0.95%  '
$ws.Range("E44").Value = '  +0.97%  '
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("E46").Value = '  -12.40%  '
$ws.Range("E47").Value = '  +2.25%  '
$ws.Range("D48").Value = "'15.34"
$ws.Range("E48").Value = '  -1.01%  '
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").Value = "'3.03"
$ws.Range("E49").Value = '  +2.84%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = "'7.29"
$ws.Range("E50").Value = '  +1.82%  '
$ws.Range("D51").Value = '2.299.16'
$ws.Range("E51").Value = '  +2.15%  '
